# ValueSet-med-addition-peri-transplant-vs.xlsx
# "update with new logo and colors" -- refresh the ValueSet metadata table:
#   - bump Version 0.1.6 -> 0.1.7
#   - Status active -> draft
#   - Date refreshed
#   - Contact reworked to the CIBMTR org contact + a new named contact (Bob Milius)
#   - new Jurisdiction row inserted between the Contact rows and Description

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- simple value updates (rows 1-11 keep their existing row numbers) ---
$ws.Range("B3").Value  = "0.1.7"
$ws.Range("B6").Value  = "draft"
$ws.Range("B8").Value  = "2024-08-23T10:17:11-05:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- insert a new "Jurisdiction" row right after the second Contact row ---
# (this pushes Description/Purpose/Copyright/Immutable down by one row,
#  from rows 12-15 to rows 13-16)
$ws.Rows.Item(12).EntireRow.Insert()

$dst = $ws.Range("A12:B12")
$src = $ws.Range("A13:B13")
$dst.Borders.LineStyle = $src.Borders.LineStyle
$dst.VerticalAlignment = $src.VerticalAlignment
$dst.WrapText = $src.WrapText

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
